$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.088.68"
$ws.Range("E2").Value = "  -1.37%  "
$ws.Range("D3").Value = "3.529.48"
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.52%  "
$ws.Range("D7").Value = "3.527.64"
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.136"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.71%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "8.01"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.31%  "
$ws.Range("E12").Value = "  -2.54%  "
$ws.Range("D13").Value = "4.129.20"
$ws.Range("E13").Value = "  +0.44%  "
$ws.Range("E14").Value = "  -4.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "30.42"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.92%  "
$ws.Range("D16").Value = "3.527.75"
$ws.Range("E16").Value = "  +0.42%  "
$ws.Range("D17").Value = "66.224.94"
$ws.Range("E17").Value = "  -1.16%  "
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.97"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.94%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "425.52"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.63%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.602"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.94%  "
$ws.Range("D25").Value = "3.675.76"
$ws.Range("E25").Value = "  +0.86%  "
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("E27").Value = "  -1.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.03"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.75%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.47"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.44%  "
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("E32").Value = "  -5.44%  "
$ws.Range("E33").Value = "  -3.88%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.34"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.38%  "
$ws.Range("D35").Value = "3.520.81"
$ws.Range("E35").Value = "  +0.14%  "
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("E37").Value = "  -2.94%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.86"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.67%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.61"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.03%  "
$ws.Range("E40").Value = "  +0.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "170.89"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0861"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.76%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.17"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.892"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.51%  "
$ws.Range("E45").Value = "  -9.96%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "45.31"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.43%  "
$ws.Range("E47").Value = "  -10.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "25.91"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -8.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.41"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.18"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.950"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.25%  "
